$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transactions")

# Fix "date import issue": the Start/End date columns (D:E) for rows
# 101-150 were left blank by the importer. Fill them in with the
# corrected date value.
$ws.Range("D101:E150").Value = 21916

# Leave the sheet in the same selection state the author ended up with:
# columns D:E fully selected (this was likely done to double-check/format
# the newly imported date values).
$ws.Range("D:E").Select() | Out-Null
